# Fixed #366 User content is lost after two generation without edition.
#
# The two `w:fldSimple` "simple field" elements (m:usercontent / m:endusercontent)
# are rewritten as "complex field" run sequences:
#   <w:r><w:fldChar w:fldCharType="begin"/></w:r>
#   <w:r><w:instrText>INSTR</w:instrText></w:r>
#   <w:r><w:fldChar w:fldCharType="separate"/></w:r>
#   <w:r><w:fldChar w:fldCharType="end"/></w:r>

$d = $word.ActiveDocument

function Expand-SimpleField($paragraphIndex, $instr) {
    $para = $d.Paragraphs($paragraphIndex)
    $pRange = $para.Range
    $pXml = $pRange.WordOpenXML

    # Grab the paragraph's own attributes (rsid*, etc.) straight out of its
    # current OOXML so the rewritten paragraph keeps them intact. Strip any
    # xmlns declarations and the synthetic w14:paraId/w14:textId attributes
    # that the WordOpenXML getter adds but which are not present in the
    # original document.
    $pAttrs = ""
    if ($pXml -match '<w:p\b([^>]*)>') {
        $pAttrs = $matches[1] -replace '\s+xmlns:\w+="[^"]*"', ''
        $pAttrs = $pAttrs -replace '\s+w14:paraId="[^"]*"', ''
        $pAttrs = $pAttrs -replace '\s+w14:textId="[^"]*"', ''
    }

    # Target range = whole paragraph minus the trailing paragraph mark, so the
    # end-of-paragraph character (and the <w:p> element itself) stays put.
    $target = $d.Range($pRange.Start, $pRange.End - 1)

    $inner = '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' + `
             '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' + `
             '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' + `
             '<w:r><w:fldChar w:fldCharType="end"/></w:r>'

    $newParaOpen = '<w:p' + $pAttrs + '>'

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData>' + `
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body>' + $newParaOpen + $inner + '</w:p></w:body></w:document>' + `
           '</pkg:xmlData></pkg:part></pkg:package>'

    $target.InsertXML($xml)
}

# Find every simple-field paragraph (<w:fldSimple w:instr="...">) together
# with its instruction text.
$targets = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $xml = $d.Paragraphs($i).Range.WordOpenXML
    if ($xml -match '<w:fldSimple w:instr="([^"]*)"') {
        $targets += , @($i, $matches[1])
    }
}

# Expand each one in turn, starting from the last paragraph so that earlier
# paragraph indices are unaffected by a later rewrite.

for ($j = $targets.Count - 1; $j -ge 0; $j--) {
    $pair = $targets[$j]
    Expand-SimpleField $pair[0] $pair[1]
}

Write-Output "Expanded $($targets.Count) simple field(s) into complex fields."
